# tracking hydroprogress and updating reduction efficiencies.
#
# - HUC 70802060502 (Mud Creek, row 3) now has hydrowork tracked:
#     Status = "In progress", Technician = "Thomas Kosacz"
# - HUC 70801030407 (Rock Creek, row 7) hydrowork is finished:
#     Status updated from "In progress" to "Complete"
# - Leave selection on G10, matching where the editor left off.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New tracking entry for Mud Creek (row 3): status + technician assigned.
$ws.Range("F3").Value = "In progress"
$ws.Range("G3").Value = "Thomas Kosacz"

# Rock Creek (row 7) hydrowork has been completed.
$ws.Range("F7").Value = "Complete"

# Leave the active selection where the author ended up.
$ws.Range("G10").Select()
